# Atualização automática de SANTO_ANTONIO_DA_PATRULHA.xlsx
#
# The edit consists of:
#   1. Renaming sheet "Paineis DARQ" -> "PAINEIS DARQ"
#   2. Renaming sheet "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#   3. Deleting the "Desarquivamentos Pendentes" sheet entirely
#      (its two header strings "PEDIDOS PENDENTES"/"%" and its five
#      dedicated cell styles disappear with it; the workbook's remaining
#      sheets/files get repacked/renumbered by the save routine.)

$wb = $excel.ActiveWorkbook

# 1. Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"

# 2. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# 3. Remove the "Desarquivamentos Pendentes" worksheet completely
$excel.DisplayAlerts = $false
$wb.Worksheets("Desarquivamentos Pendentes").Delete() | Out-Null
